$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "69.068.79"
Set-TextValue "E2" "  -2.11%  "
Set-TextValue "D3" "3.676.79"
Set-TextValue "E3" "  -2.93%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "679.33"
Set-TextValue "E5" "  -3.92%  "
Set-TextValue "D6" "161.44"
Set-TextValue "E6" "  -4.78%  "
Set-TextValue "D7" "3.676.04"
Set-TextValue "E7" "  -2.90%  "
Set-TextValue "E8" "  +0.20%  "
Set-TextValue "D9" "0.494"
Set-TextValue "E9" "  -4.96%  "
Set-TextValue "E10" "  -7.47%  "
Set-TextValue "D11" "7.24"
Set-TextValue "E11" "  -1.36%  "
Set-TextValue "D12" "0.450"
Set-TextValue "E12" "  -1.34%  "
Set-TextValue "E13" "  -7.15%  "
Set-TextValue "D14" "33.23"
Set-TextValue "E14" "  -7.91%  "
Set-TextValue "D15" "4.298.82"
Set-TextValue "E15" "  -2.67%  "
Set-TextValue "D16" "3.675.79"
Set-TextValue "E16" "  -3.93%  "
Set-TextValue "D17" "69.147.38"
Set-TextValue "E17" "  -1.90%  "
Set-TextValue "E18" "  -1.70%  "
Set-TextValue "E19" "  -6.29%  "
Set-TextValue "D20" "6.61"
Set-TextValue "E20" "  -7.19%  "
Set-TextValue "D21" "483.03"
Set-TextValue "E21" "  -1.76%  "
Set-TextValue "D22" "9.76"
Set-TextValue "E22" "  -7.94%  "
Set-TextValue "D23" "0.661"
Set-TextValue "E23" "  -8.87%  "
Set-TextValue "D24" "79.42"
Set-TextValue "E24" "  -6.36%  "
Set-TextValue "D25" "3.821.82"
Set-TextValue "E25" "  -2.74%  "
Set-TextValue "D26" "0.0000127"
Set-TextValue "E26" "  -12.06%  "
Set-TextValue "B27" "InternetComputer(DFINITY)"
Set-TextValue "C27" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D27" "11.53"
Set-TextValue "E27" "  -4.41%  "
Set-TextValue "B28" "Dai"
Set-TextValue "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D28" "0.999"
Set-TextValue "E28" "  -0.05%  "
Set-TextValue "D29" "9.47"
Set-TextValue "E29" "  -9.35%  "
Set-TextValue "D30" "1.83"
Set-TextValue "E30" "  -10.64%  "
Set-TextValue "E31" "  -11.44%  "
Set-TextValue "E32" "  -5.42%  "
Set-TextValue "D33" "6.70"
Set-TextValue "E33" "  -8.26%  "
Set-TextValue "E34" "  +0.12%  "
Set-TextValue "D35" "26.82"
Set-TextValue "E35" "  -7.74%  "
Set-TextValue "E36" "  -6.71%  "
Set-TextValue "D37" "3.642.18"
Set-TextValue "E37" "  -2.96%  "
Set-TextValue "D38" "8.48"
Set-TextValue "E38" "  -5.87%  "
Set-TextValue "D39" "6.03"
Set-TextValue "E39" "  +2.00%  "
Set-TextValue "D40" "0.0936"
Set-TextValue "E40" "  -7.24%  "
Set-TextValue "E41" "  -0.02%  "
Set-TextValue "D42" "2.16"
Set-TextValue "E42" "  -6.55%  "
Set-TextValue "E43" "  +0.13%  "
Set-TextValue "D44" "0.956"
Set-TextValue "E44" "  -7.84%  "
Set-TextValue "D45" "157.27"
Set-TextValue "E45" "  -4.13%  "
Set-TextValue "D46" "47.97"
Set-TextValue "E46" "  -1.85%  "
Set-TextValue "D47" "2.79"
Set-TextValue "E47" "  -14.85%  "
Set-TextValue "E48" "  -11.49%  "
Set-TextValue "D49" "389.84"
Set-TextValue "E49" "  -7.31%  "
Set-TextValue "E50" "  -5.38%  "
Set-TextValue "B51" "InjectiveProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "28.15"
Set-TextValue "E51" "  +1.39%  "
